# Auto-generated edit script applying numeric value updates per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 9370.333000000001
$ws.Range("J69").Value = 9370.333000000001
$ws.Range("L69").Value = 28110.999
$ws.Range("N69").Value = -29858.999
$ws.Range("H72").Value = 9370.333000000001
$ws.Range("J72").Value = 9370.333000000001
$ws.Range("L72").Value = 84332.997
$ws.Range("N72").Value = -93068.997
$ws.Range("H74").Value = 9437.875
$ws.Range("J74").Value = 9208.666999999999
$ws.Range("L74").Value = 9208.666999999999
$ws.Range("N74").Value = -11080.667
$ws.Range("H77").Value = 9437.875
$ws.Range("J77").Value = 9208.666999999999
$ws.Range("L77").Value = 46043.335
$ws.Range("N77").Value = -55403.335
$ws.Range("H112").Value = 1415.2142
$ws.Range("H129").Value = 1493
$ws.Range("J129").Value = 3500
$ws.Range("L129").Value = 10500
$ws.Range("N129").Value = -20500
$ws.Range("H132").Value = 2011.1086
$ws.Range("I132").Value = 1443.3954
$ws.Range("J132").Value = 10148.333
$ws.Range("K132").Value = 4330.1862
$ws.Range("L132").Value = 30444.999
$ws.Range("M132").Value = -1800.1862
$ws.Range("N132").Value = -35504.999
$ws.Range("H138").Value = 3471.0981
$ws.Range("I138").Value = 2557.4092
$ws.Range("J138").Value = 4164.241
$ws.Range("K138").Value = 7672.2276
$ws.Range("L138").Value = 12492.723
$ws.Range("M138").Value = -2532.2276
$ws.Range("N138").Value = -22772.723
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 16197.6
$ws.Range("J28").Value = 17999.5
$ws.Range("L28").Value = 17999.5
$ws.Range("N28").Value = -18383.5
$ws.Range("H32").Value = 1871
$ws.Range("I32").Value = 1975.4048
$ws.Range("K32").Value = 1975.4048
$ws.Range("M32").Value = -1688.4048
$ws.Range("H45").Value = 3866
$ws.Range("I45").Value = 2928.0527
$ws.Range("J45").Value = 5846.1113
$ws.Range("K45").Value = 2928.0527
$ws.Range("L45").Value = 5846.1113
$ws.Range("M45").Value = -2551.0527
$ws.Range("N45").Value = -6600.1113
$ws.Range("H61").Value = 6228.7085
$ws.Range("I61").Value = 3962.75
$ws.Range("K61").Value = 3962.75
$ws.Range("M61").Value = -3750.75
$ws.Range("H99").Value = 16197.6
$ws.Range("J99").Value = 17999.5
$ws.Range("L99").Value = 17999.5
$ws.Range("N99").Value = -23989.5
$ws.Range("H122").Value = 3343.65
$ws.Range("I122").Value = 2901.7812
$ws.Range("K122").Value = 8705.3436
$ws.Range("M122").Value = -6255.3436
$ws.Range("H132").Value = 7327.2188
$ws.Range("I132").Value = 5471.364
$ws.Range("J132").Value = 11410.1
$ws.Range("K132").Value = 16414.092
$ws.Range("L132").Value = 34230.3
$ws.Range("M132").Value = -13884.092
$ws.Range("N132").Value = -39290.3
$ws.Range("H136").Value = 6228.7085
$ws.Range("I136").Value = 3962.75
$ws.Range("K136").Value = 11888.25
$ws.Range("M136").Value = -9338.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12665.667
$ws.Range("I86").Value = 9495
$ws.Range("J86").Value = 19007
$ws.Range("K86").Value = 9495
$ws.Range("L86").Value = 19007
$ws.Range("M86").Value = -8372
$ws.Range("N86").Value = -21253
$ws.Range("H89").Value = 12665.667
$ws.Range("I89").Value = 9495
$ws.Range("J89").Value = 19007
$ws.Range("K89").Value = 47475
$ws.Range("L89").Value = 95035
$ws.Range("M89").Value = -41859
$ws.Range("N89").Value = -106267
$ws.Range("H134").Value = 2862.95
$ws.Range("I134").Value = 1958.0555
$ws.Range("K134").Value = 5874.166499999999
$ws.Range("M134").Value = -3339.166499999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5500.227
$ws.Range("I58").Value = 2591.5
$ws.Range("K58").Value = 2591.5
$ws.Range("M58").Value = -2388.5
$ws.Range("H99").Value = 3658.1428
$ws.Range("I99").Value = 3434.5
$ws.Range("K99").Value = 3434.5
$ws.Range("M99").Value = -1936.5
$ws.Range("H126").Value = 3658.1428
$ws.Range("I126").Value = 3434.5
$ws.Range("K126").Value = 10303.5
$ws.Range("M126").Value = -7833.5
$ws.Range("H132").Value = 2074.6775
$ws.Range("I132").Value = 1450.1111
$ws.Range("K132").Value = 4350.3333
$ws.Range("M132").Value = -1820.3333
$ws.Range("H134").Value = 2046.5454
$ws.Range("I134").Value = 1222.8206
$ws.Range("J134").Value = 8471.6
$ws.Range("K134").Value = 3668.4618
$ws.Range("L134").Value = 25414.8
$ws.Range("M134").Value = -1133.4618
$ws.Range("N134").Value = -30484.8
$ws.Range("H136").Value = 5500.227
$ws.Range("I136").Value = 2591.5
$ws.Range("K136").Value = 7774.5
$ws.Range("M136").Value = -5224.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 84348.58
$ws.Range("I2").Value = 93.375
$ws.Range("J2").Value = 145625.1
$ws.Range("K2").Value = 560.25
$ws.Range("L2").Value = 873750.6000000001
$ws.Range("M2").Value = -447.25
$ws.Range("N2").Value = -873976.6000000001
$ws.Range("H3").Value = 1230.8572
$ws.Range("I3").Value = 1230.8572
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3692.5716
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -3580.5716
$ws.Range("N3").ClearContents()
$ws.Range("H11").Value = 35766.758
$ws.Range("I11").Value = 54090.74
$ws.Range("J11").Value = 951.2
$ws.Range("K11").Value = 162272.22
$ws.Range("L11").Value = 2853.6
$ws.Range("M11").Value = -162132.22
$ws.Range("N11").Value = -3133.6
$ws.Range("H17").Value = 270.66666
$ws.Range("I17").Value = 162
$ws.Range("K17").Value = 486
$ws.Range("M17").Value = -317
$ws.Range("H133").Value = 21515
$ws.Range("I133").Value = 21515
$ws.Range("K133").Value = 64545
$ws.Range("M133").Value = -59485
$ws.Range("H136").Value = 2315.8
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 2272
$ws.Range("I137").Value = 699.6667
$ws.Range("J137").Value = 3844.3333
$ws.Range("K137").Value = 2099.0001
$ws.Range("L137").Value = 11532.9999
$ws.Range("M137").Value = 3000.9999
$ws.Range("N137").Value = -21732.9999
$ws.Range("H138").Value = 8750
$ws.Range("J138").Value = 15000
$ws.Range("L138").Value = 45000
$ws.Range("N138").Value = -55280
$ws.Range("H139").Value = 6594.8335
$ws.Range("J139").Value = 7231.8125
$ws.Range("L139").Value = 21695.4375
$ws.Range("N139").Value = -31975.4375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6790.727
$ws.Range("I80").Value = 4615.5
$ws.Range("J80").Value = 9401
$ws.Range("K80").Value = 4615.5
$ws.Range("L80").Value = 9401
$ws.Range("M80").Value = -3617.5
$ws.Range("N80").Value = -11397
$ws.Range("H83").Value = 6790.727
$ws.Range("I83").Value = 4615.5
$ws.Range("J83").Value = 9401
$ws.Range("K83").Value = 23077.5
$ws.Range("L83").Value = 47005
$ws.Range("M83").Value = -18085.5
$ws.Range("N83").Value = -56989
$ws.Range("H97").Value = 1317.65
$ws.Range("I97").Value = 956.63635
$ws.Range("K97").Value = 956.63635
$ws.Range("M97").Value = -460.63635
$ws.Range("H102").Value = 2863.5454
$ws.Range("I102").Value = 2062.375
$ws.Range("K102").Value = 2062.375
$ws.Range("M102").Value = -440.375
$ws.Range("H132").Value = 10175
$ws.Range("I132").Value = 7505.5
$ws.Range("K132").Value = 22516.5
$ws.Range("M132").Value = -19986.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2632263.5
$ws.Range("I55").Value = 3333607.2
$ws.Range("K55").Value = 3333607.2
$ws.Range("M55").Value = -3333434.2
$ws.Range("H93").Value = 24761.455
$ws.Range("I93").Value = 19370.143
$ws.Range("J93").Value = 34196.25
$ws.Range("K93").Value = 19370.143
$ws.Range("L93").Value = 34196.25
$ws.Range("M93").Value = -18122.143
$ws.Range("N93").Value = -36692.25
$ws.Range("H122").Value = 6887.5
$ws.Range("I122").Value = 4849.5
$ws.Range("K122").Value = 14548.5
$ws.Range("M122").Value = -12098.5
$ws.Range("H132").Value = 5791.7036
$ws.Range("I132").Value = 4039.1333
$ws.Range("K132").Value = 12117.3999
$ws.Range("M132").Value = -9587.3999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14433.777
$ws.Range("I81").Value = 9699.75
$ws.Range("K81").Value = 19399.5
$ws.Range("M81").Value = -18338.5
$ws.Range("H84").Value = 14433.777
$ws.Range("I84").Value = 9699.75
$ws.Range("K84").Value = 96997.5
$ws.Range("M84").Value = -91693.5
$ws.Range("H132").Value = 6036.659
$ws.Range("I132").Value = 2688.861
$ws.Range("J132").Value = 21101.75
$ws.Range("K132").Value = 8066.583
$ws.Range("L132").Value = 63305.25
$ws.Range("M132").Value = -5536.583
$ws.Range("N132").Value = -68365.25
